# Pavullo.xlsx report update ("aggiornato a 2/3, aggiornati i report")
# - fills in a previously missing day (A=44235) between old rows 92 and 93
# - shifts the following rows down by one
# - appends a brand-new day (A=44257) at the bottom
# - refreshes the rolling 7-day sums ("somma mobile 7gg.") affected by the
#   newly inserted day

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Shift rows 93..113 down to 94..114 (bottom-up so we never read a cell
#    after it has already been overwritten).
for ($r = 113; $r -ge 93; $r--) {
    $dst = $r + 1
    $ws.Cells.Item($dst, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($dst, 2).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($dst, 3).Value = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($dst, 4).Value = $ws.Cells.Item($r, 4).Value2
}

# 2) The 7-day rolling sum/rate for what is now row 112 (previously row 111,
#    which had no C/D yet) is now computable.
$ws.Cells.Item(112, 3).Value = 31
$ws.Cells.Item(112, 4).Value = 172.4329736344421

# 3) Write the newly-available day into the row that opened up at 93.
$ws.Cells.Item(93, 1).Value = 44235
$ws.Cells.Item(93, 2).Value = 5
$ws.Cells.Item(93, 3).Value = 26
$ws.Cells.Item(93, 4).Value = 144.621203693403

# 4) Recompute the rolling 7-day window for the days whose window now
#    includes the newly-filled-in day.
$ws.Cells.Item(90, 3).Value = 26
$ws.Cells.Item(90, 4).Value = 144.621203693403

$ws.Cells.Item(91, 3).Value = 23
$ws.Cells.Item(91, 4).Value = 127.9341417287796

$ws.Cells.Item(92, 3).Value = 24
$ws.Cells.Item(92, 4).Value = 133.4964957169874

# 5) Append the new final day at the bottom of the table.
$ws.Cells.Item(115, 1).Value = 44257
$ws.Cells.Item(115, 2).Value = 6

# Rows beyond the sheet's original extent (114, 115) don't inherit the date
# column's style when written through Cells.Item, so copy it explicitly from
# an existing date cell (format only - this doesn't disturb the values just
# written above).
$ws.Range("A92").Copy()
$ws.Range("A114:A115").PasteSpecial(-4122)
